# Update the subtitle on slide 1 to reference the new MPS version
# (2021.1.x -> 2022.2.x), per the "Adapt documentation to 2022.2 migration"
# commit.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -like "MPS 2021.1.x*") {
            $sh.TextFrame.TextRange.Text = "MPS 2022.2.x + mbeddr.platform / MPS Extensions + KernelF"
        }
    }
}
